# Auto-generated edit script applying cell-level numeric updates
# to the 'Profit' analysis columns (H-N) across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H129").Value = 885.4
$ws.Range("J129").Value = 878.4039
$ws.Range("L129").Value = 2635.2117
$ws.Range("N129").Value = -12635.2117
$ws.Range("H131").Value = 2549.4167
$ws.Range("J131").Value = 3188.6667
$ws.Range("L131").Value = 9566.000100000001
$ws.Range("N131").Value = -19646.0001
$ws.Range("H132").Value = 791.0789
$ws.Range("I132").Value = 734.32355
$ws.Range("J132").Value = 1273.5
$ws.Range("K132").Value = 2202.97065
$ws.Range("L132").Value = 3820.5
$ws.Range("M132").Value = 327.0293500000002
$ws.Range("N132").Value = -8880.5
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 2039.3704
$ws.Range("I137").Value = 1163.5
$ws.Range("J137").Value = 2191.6956
$ws.Range("K137").Value = 3490.5
$ws.Range("L137").Value = 6575.0868
$ws.Range("M137").Value = -940.5
$ws.Range("N137").Value = -11675.0868
$ws.Range("H138").Value = 1768.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1924.4022
$ws.Range("I32").Value = 1421.3827
$ws.Range("J32").Value = 5628.4546
$ws.Range("K32").Value = 1421.3827
$ws.Range("L32").Value = 5628.4546
$ws.Range("M32").Value = -1134.3827
$ws.Range("N32").Value = -6202.4546
$ws.Range("H45").Value = 3752790.8
$ws.Range("I45").Value = 12859883
$ws.Range("J45").Value = 2811.7058
$ws.Range("K45").Value = 12859883
$ws.Range("L45").Value = 2811.7058
$ws.Range("M45").Value = -12859506
$ws.Range("N45").Value = -3565.7058
$ws.Range("H61").Value = 3008.087
$ws.Range("I61").Value = 2023.7333
$ws.Range("K61").Value = 2023.7333
$ws.Range("M61").Value = -1811.7333
$ws.Range("H74").Value = 1599.1875
$ws.Range("I74").Value = 606
$ws.Range("K74").Value = 606
$ws.Range("M74").Value = 268
$ws.Range("H77").Value = 1599.1875
$ws.Range("I77").Value = 606
$ws.Range("K77").Value = 3030
$ws.Range("M77").Value = 1338
$ws.Range("H136").Value = 3008.087
$ws.Range("I136").Value = 2023.7333
$ws.Range("K136").Value = 6071.199900000001
$ws.Range("M136").Value = -3521.199900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 35000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 35000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -35620
$ws.Range("H99").Value = 1157.4615
$ws.Range("I99").Value = 1210.7
$ws.Range("K99").Value = 1210.7
$ws.Range("M99").Value = 287.3
$ws.Range("H105").Value = 2519.9
$ws.Range("I105").Value = 2494.6316
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2494.6316
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -747.6316000000002
$ws.Range("N105").Value = -6494
$ws.Range("H107").Value = 5980
$ws.Range("I107").Value = 5980
$ws.Range("K107").Value = 5980
$ws.Range("M107").Value = -4060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 882.8
$ws.Range("I16").Value = 857.55554
$ws.Range("J16").Value = 1110
$ws.Range("K16").Value = 857.55554
$ws.Range("L16").Value = 1110
$ws.Range("M16").Value = -570.55554
$ws.Range("N16").Value = -1684
$ws.Range("H31").Value = 1441.7878
$ws.Range("I31").Value = 1024.7142
$ws.Range("J31").Value = 1554.0769
$ws.Range("K31").Value = 1024.7142
$ws.Range("L31").Value = 1554.0769
$ws.Range("M31").Value = -729.7141999999999
$ws.Range("N31").Value = -2144.0769
$ws.Range("H34").Value = 1441.7878
$ws.Range("I34").Value = 1024.7142
$ws.Range("J34").Value = 1554.0769
$ws.Range("K34").Value = 1024.7142
$ws.Range("L34").Value = 1554.0769
$ws.Range("M34").Value = -822.7141999999999
$ws.Range("N34").Value = -1958.0769
$ws.Range("H58").Value = 2416853
$ws.Range("I58").Value = 4349320
$ws.Range("K58").Value = 4349320
$ws.Range("M58").Value = -4349117
$ws.Range("H70").Value = 28833.334
$ws.Range("J70").Value = 28833.334
$ws.Range("L70").Value = 28833.334
$ws.Range("N70").Value = -29463.334
$ws.Range("H73").Value = 28833.334
$ws.Range("J73").Value = 28833.334
$ws.Range("L73").Value = 28833.334
$ws.Range("N73").Value = -31017.334
$ws.Range("H107").Value = 819.3077
$ws.Range("I107").Value = 762.5833
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 762.5833
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1157.4167
$ws.Range("N107").Value = -5340
$ws.Range("H113").Value = 882.8
$ws.Range("I113").Value = 857.55554
$ws.Range("J113").Value = 1110
$ws.Range("K113").Value = 857.55554
$ws.Range("L113").Value = 1110
$ws.Range("M113").Value = 1312.44446
$ws.Range("N113").Value = -5450
$ws.Range("H136").Value = 2416853
$ws.Range("I136").Value = 4349320
$ws.Range("K136").Value = 13047960
$ws.Range("M136").Value = -13045410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2249.5
$ws.Range("H89").Value = 2249.5
$ws.Range("H117").Value = 28571742
$ws.Range("J117").Value = 71428660
$ws.Range("L117").Value = 214285980
$ws.Range("N117").Value = -214292864
$ws.Range("H131").Value = 9275254
$ws.Range("J131").Value = 17591.674
$ws.Range("L131").Value = 52775.022
$ws.Range("N131").Value = -62855.022
$ws.Range("H134").Value = 4725.154
$ws.Range("I134").Value = 4404.5
$ws.Range("K134").Value = 13213.5
$ws.Range("M134").Value = -8143.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19566.666
$ws.Range("J46").Value = 19566.666
$ws.Range("L46").Value = 19566.666
$ws.Range("N46").Value = -19878.666
$ws.Range("H97").Value = 1382.125
$ws.Range("I97").Value = 991.1111
$ws.Range("J97").Value = 1884.8572
$ws.Range("K97").Value = 991.1111
$ws.Range("L97").Value = 1884.8572
$ws.Range("M97").Value = -495.1111
$ws.Range("N97").Value = -2876.8572
$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 4500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -2254
$ws.Range("H104").Value = 49997
$ws.Range("J104").Value = 49997
$ws.Range("L104").Value = 49997
$ws.Range("N104").Value = -56985
$ws.Range("H122").Value = 2416.2942
$ws.Range("I122").Value = 2129.75
$ws.Range("J122").Value = 2671
$ws.Range("K122").Value = 6389.25
$ws.Range("L122").Value = 8013
$ws.Range("M122").Value = -3939.25
$ws.Range("N122").Value = -12913
$ws.Range("H132").Value = 1481713.6
$ws.Range("I132").Value = 2026202.4
$ws.Range("K132").Value = 6078607.199999999
$ws.Range("M132").Value = -6076077.199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3683.5386
$ws.Range("I61").Value = 3375.75
$ws.Range("J61").Value = 3820.3333
$ws.Range("K61").Value = 3375.75
$ws.Range("L61").Value = 3820.3333
$ws.Range("M61").Value = -3173.75
$ws.Range("N61").Value = -4224.3333
$ws.Range("H93").Value = 1395.6
$ws.Range("I93").Value = 992.6667
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 992.6667
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 255.3333
$ws.Range("N93").Value = -4496
$ws.Range("H113").Value = 3683.5386
$ws.Range("I113").Value = 3375.75
$ws.Range("J113").Value = 3820.3333
$ws.Range("K113").Value = 3375.75
$ws.Range("L113").Value = 3820.3333
$ws.Range("M113").Value = -1205.75
$ws.Range("N113").Value = -8160.3333
$ws.Range("H132").Value = 4186.7085
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4186.7085
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12560.1255
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -17620.1255
$ws.Range("H136").Value = 3673.5557
$ws.Range("I136").Value = 2347.6843
$ws.Range("J136").Value = 6822.5
$ws.Range("K136").Value = 7043.0529
$ws.Range("L136").Value = 20467.5
$ws.Range("M136").Value = -4493.0529
$ws.Range("N136").Value = -25567.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 75005
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H132").Value = 1594.5897
$ws.Range("I132").Value = 1318.3704
$ws.Range("J132").Value = 2216.0833
$ws.Range("K132").Value = 3955.1112
$ws.Range("L132").Value = 6648.249899999999
$ws.Range("M132").Value = -1425.1112
$ws.Range("N132").Value = -11708.2499
$ws.Range("H136").Value = 12348967
$ws.Range("I136").Value = 19844328
$ws.Range("K136").Value = 59532984
$ws.Range("M136").Value = -59530434
